$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows among the "Source:" footer block, pushing the
# existing rows down (matching the target layout: 47,48(blank),49,50(blank),
# 51,52(blank),53,...,56,57).
$ws.Rows(48).Insert()
$ws.Rows(50).Insert()
$ws.Rows(52).Insert()

# Make sure the newly inserted blank rows use the "source" (italic) style,
# same as the rest of this footer block.
$ws.Range("A48").Style = "source"
$ws.Range("A50").Style = "source"
$ws.Range("A52").Style = "source"

# The hyperlinked URL cell (now at A51) loses its hyperlink styling and
# becomes plain "source" text like its neighbours.
$ws.Range("A51").Style = "source"

# Remove the hyperlink itself.
$ws.Hyperlinks.Delete()

# The long citation text in the final row is replaced by the short "PCBS"
# text (matching the cell above it).
$ws.Range("A57").Value = "PCBS"
